# feat: add 2022-Q3 data
#
# The workbook has a summary sheet ("总计") and one quarterly sheet
# ("2022-Q2"). This change introduces a new "2022-Q3" quarterly sheet
# (with fresh fund data) while keeping the old "2022-Q2" sheet (with its
# original data) as a separate tab placed after it. The summary sheet
# gets an extra row recording the "2022-Q2" entry that used to be on row 2,
# and row 2 itself is relabelled to "2022-Q3".

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the existing "2022-Q2" sheet -----------------------
# The duplicate keeps the original data/styling untouched and will stay
# named "2022-Q2"; the original sheet will be repurposed below to hold
# the new "2022-Q3" figures.
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Copy($null, $oldQ2)
$dupQ2 = $wb.Worksheets.Item($oldQ2.Index + 1)

# --- 2. Turn the original sheet into "2022-Q3" with new figures ------
# (D2:G2 are text cells, like the originals - prefix with an apostrophe so
# the numeric-looking strings aren't coerced to numbers, then strip the
# resulting quote-prefix formatting by repainting with B2's plain format.)
$oldQ2.Range("D2").Value = "'1.12"
$oldQ2.Range("E2").Value = "'90.06"
$oldQ2.Range("F2").Value = "'2.57"
$oldQ2.Range("G2").Value = "'0.0288"
$oldQ2.Range("H2").Value = 9
$oldQ2.Range("B2").Copy()
$oldQ2.Range("D2:G2").PasteSpecial(-4122)

# --- 3. Rename the sheets (order: 总计, 2022-Q3, 2022-Q2) ------------
$oldQ2.Name = "2022-Q3"
$dupQ2.Name = "2022-Q2"

# --- 4. Update the "总计" summary sheet -------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("B2").Value = "2022-Q3"

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.03

# Match A2's cell formatting (bold/bordered) on the new A3 cell.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

# Keep the originally-active "总计" tab selected (sheet additions/copies
# above shift Excel's active sheet as a side effect).
$summary.Activate()
